# "Generate Report for Handoff" - updates localization-status report
# to reflect that e2e\b.md is now ready for handoff (new xliff files
# generated) while a.md stays as previously handed back.

$wb = $excel.ActiveWorkbook

$statusReady  = "Ready for handoff"
$newHandoffDateTime_zhcn = "2016-08-22 10:17:36"
$newHandoffDateTime_dede = "2016-08-22 10:17:41"
$newHandoffFile_zhcn = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$newHandoffFile_dede = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e6f4acdd34021ed04077c9a0ace7bd8ce895d19e/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7370578cfcd52181f7703f24cb50baecc28cc6b6/e2e/b.md."

# ---------------------------------------------------------------
# Overview sheet: row 3 is the b.md summary row
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = $newHandoffDateTime_dede

# ---------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md detail row
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("G3").Value = $newHandoffFile_zhcn
$wsZhCn.Range("H3").Value = $newHandoffDateTime_zhcn
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------
# de-de sheet: row 3 is the b.md detail row
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("G3").Value = $newHandoffFile_dede
$wsDeDe.Range("H3").Value = $newHandoffDateTime_dede
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
